$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with uniform run formatting) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Crime table updates (rows 14-29) ---
$ws.Range("N14").Value = -95.454545454545
$ws.Range("G15").Value = 2
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 5
$ws.Range("C14").Copy()
$ws.Range("D16").PasteSpecial(-4104)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4104)
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 36
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = 9.090909090909
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -48.571428571428
$ws.Range("N16").Value = -85.826771653543
$ws.Range("C17").Value = 12
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 3
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 5.263157894736
$ws.Range("I17").Value = 73
$ws.Range("J17").Value = 82
$ws.Range("K17").Value = -10.975609756097
$ws.Range("L17").Value = 23.728813559322
$ws.Range("M17").Value = 55.31914893617
$ws.Range("N17").Value = -69.709543568464
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4104)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = -23.529411764705
$ws.Range("L18").Value = -51.851851851851
$ws.Range("M18").Value = -13.333333333333
$ws.Range("N18").Value = -90.510948905109
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -85.714285714285
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -44.827586206896
$ws.Range("I19").Value = 92
$ws.Range("J19").Value = 126
$ws.Range("K19").Value = -26.984126984127
$ws.Range("L19").Value = -4.166666666666
$ws.Range("M19").Value = 84
$ws.Range("N19").Value = -20
$ws.Range("C20").Value = 1
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 1
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 36
$ws.Range("J20").Value = 22
$ws.Range("K20").Value = 63.636363636363
$ws.Range("L20").Value = 200
$ws.Range("M20").Value = 56.521739130434
$ws.Range("N20").Value = -65.04854368932
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 35.714285714285
$ws.Range("F21").Value = 56
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = -16.417910447761
$ws.Range("I21").Value = 265
$ws.Range("J21").Value = 304
$ws.Range("K21").Value = -12.828947368421
$ws.Range("L21").Value = 2.316602316602
$ws.Range("M21").Value = 14.718614718614
$ws.Range("N21").Value = -74.321705426356
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4104)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4104)
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 11
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = 10
$ws.Range("L22").Value = 175
$ws.Range("M22").Value = 266.666666666667
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -66.666666666666
$ws.Range("F24").Value = 46
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = -59.29203539823
$ws.Range("I24").Value = 233
$ws.Range("J24").Value = 368
$ws.Range("K24").Value = -36.684782608695
$ws.Range("L24").Value = 12.01923076923
$ws.Range("M24").Value = 121.904761904762
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = -29.629629629629
$ws.Range("I25").Value = 101
$ws.Range("J25").Value = 133
$ws.Range("K25").Value = -24.060150375939
$ws.Range("L25").Value = 23.170731707317
$ws.Range("M25").Value = -33.112582781457
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4104)
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4104)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C28").Value = 1
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 3
$ws.Range("K28").Value = -25
$ws.Range("L28").Value = -25
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = -92.5
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("C29").Value = 1
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = -25
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = -92.307692307692
